$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.06939234730424
$ws.Range("D2").Value = 1.072038659540243
$ws.Range("E2").Value = 1.073297678660724
$ws.Range("F2").Value = 1.082757872379919
$ws.Range("I2").Value = 1.046449540882804
$ws.Range("J2").Value = 1.074327305605132
$ws.Range("K2").Value = 1.074734235453945
$ws.Range("L2").Value = 1.075989911868544
$ws.Range("M2").Value = 1.08542525720971
$ws.Range("N2").Value = 1.075852973606156
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.071334476799323
$ws.Range("D3").Value = 1.073899115812292
$ws.Range("E3").Value = 1.075088726705865
$ws.Range("F3").Value = 1.084741177176355
$ws.Range("I3").Value = 1.046949866511578
$ws.Range("J3").Value = 1.075921909107392
$ws.Range("K3").Value = 1.076408787342939
$ws.Range("L3").Value = 1.077595473086181
$ws.Range("M3").Value = 1.087224443993265
$ws.Range("N3").Value = 1.077449841628292
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.072586654219481
$ws.Range("D4").Value = 1.075098774253811
$ws.Range("E4").Value = 1.076243301880598
$ws.Range("F4").Value = 1.086020358652081
$ws.Range("I4").Value = 1.047270040749313
$ws.Range("J4").Value = 1.076948969065277
$ws.Range("K4").Value = 1.077487723750012
$ws.Range("L4").Value = 1.078629581626147
$ws.Range("M4").Value = 1.088384077969757
$ws.Range("N4").Value = 1.078478360129123
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.073112013468997
$ws.Range("D5").Value = 1.075602131529832
$ws.Range("E5").Value = 1.076727664906222
$ws.Range("F5").Value = 1.086557156710355
$ws.Range("I5").Value = 1.047403794706232
$ws.Range("J5").Value = 1.077379626532699
$ws.Range("K5").Value = 1.077940225301535
$ws.Range("L5").Value = 1.079063192836862
$ws.Range("M5").Value = 1.088870518439803
$ws.Range("N5").Value = 1.078909629179546
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.073200162347302
$ws.Range("D6").Value = 1.075686590644468
$ws.Range("E6").Value = 1.076808932378736
$ws.Range("F6").Value = 1.086647231273349
$ws.Range("I6").Value = 1.047426203110089
$ws.Range("J6").Value = 1.077451870797037
$ws.Range("K6").Value = 1.078016139367563
$ws.Range("L6").Value = 1.079135932507028
$ws.Range("M6").Value = 1.088952131846363
$ws.Range("N6").Value = 1.078981976039026
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.072593678213438
$ws.Range("D7").Value = 1.075105503956114
$ws.Range("E7").Value = 1.076249777943125
$ws.Range("F7").Value = 1.086027535152963
$ws.Range("I7").Value = 1.04727183129552
$ws.Range("J7").Value = 1.076954727899651
$ws.Range("K7").Value = 1.077493774326843
$ws.Range("L7").Value = 1.078635379965138
$ws.Range("M7").Value = 1.08839058198428
$ws.Range("N7").Value = 1.078484127141703
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.070049647230855
$ws.Range("D8").Value = 1.072668288413038
$ws.Range("E8").Value = 1.07390388609409
$ws.Range("F8").Value = 1.083429012911838
$ws.Range("I8").Value = 1.04661937189569
$ws.Range("J8").Value = 1.074867205924667
$ws.Range("K8").Value = 1.075301125168268
$ws.Range("L8").Value = 1.076533524367942
$ws.Range("M8").Value = 1.086034258376845
$ws.Range("N8").Value = 1.076393640646071
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.065531128678239
$ws.Range("D9").Value = 1.068340592487433
$ws.Range("E9").Value = 1.069735809464019
$ws.Range("F9").Value = 1.078817242197701
$ws.Range("I9").Value = 1.045441980757522
$ws.Range("J9").Value = 1.071151408966262
$ws.Range("K9").Value = 1.07140115884313
$ws.Range("L9").Value = 1.072792126366318
$ws.Range("M9").Value = 1.081846192182141
$ws.Range("N9").Value = 1.072672566829778
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.062493336651469
$ws.Range("D10").Value = 1.065431876019676
$ws.Range("E10").Value = 1.066932641567722
$ws.Range("F10").Value = 1.075719150100937
$ws.Range("I10").Value = 1.044637975996189
$ws.Range("J10").Value = 1.068647856686165
$ws.Range("K10").Value = 1.068775511838252
$ws.Range("L10").Value = 1.070271243796346
$ws.Range("M10").Value = 1.079028591400124
$ws.Range("N10").Value = 1.070165459218282
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.06117156612564
$ws.Range("D11").Value = 1.064166463013171
$ws.Range("E11").Value = 1.065712730321789
$ws.Range("F11").Value = 1.074371713685094
$ws.Range("I11").Value = 1.044285202006269
$ws.Range("J11").Value = 1.067557255895234
$ws.Range("K11").Value = 1.067632196163228
$ws.Range("L11").Value = 1.069173069446698
$ws.Range("M11").Value = 1.07780216230681
$ws.Range("N11").Value = 1.069073309649134
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.06067961292399
$ws.Range("D12").Value = 1.063695515824381
$ws.Range("E12").Value = 1.06525865451957
$ws.Range("F12").Value = 1.073870293953782
$ws.Range("I12").Value = 1.044153460467514
$ws.Range("J12").Value = 1.067151150050018
$ws.Range("K12").Value = 1.067206531868189
$ws.Range("L12").Value = 1.068764139423415
$ws.Range("M12").Value = 1.077345624669521
$ws.Range("N12").Value = 1.068666627087037
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.060785183774548
$ws.Range("D13").Value = 1.063796577506089
$ws.Range("E13").Value = 1.065356098584321
$ws.Range("F13").Value = 1.073977892376076
$ws.Range("I13").Value = 1.044181751553841
$ws.Range("J13").Value = 1.067238307153426
$ws.Range("K13").Value = 1.067297883333151
$ws.Range("L13").Value = 1.068851902830826
$ws.Range("M13").Value = 1.077443598661481
$ws.Range("N13").Value = 1.068753907963529
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.061130921462196
$ws.Range("D14").Value = 1.064127553225331
$ws.Range("E14").Value = 1.06567521573825
$ws.Range("F14").Value = 1.074330285149197
$ws.Range("I14").Value = 1.044274326657213
$ws.Range("J14").Value = 1.067523707769457
$ws.Range("K14").Value = 1.067597030881411
$ws.Range("L14").Value = 1.069139288101726
$ws.Range("M14").Value = 1.077764445062813
$ws.Range("N14").Value = 1.069039713881171
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.061343809956546
$ws.Range("D15").Value = 1.064331356032319
$ws.Range("E15").Value = 1.06587170812578
$ws.Range("F15").Value = 1.074547282934236
$ws.Range("I15").Value = 1.044331271458676
$ws.Range("J15").Value = 1.067699418173917
$ws.Range("K15").Value = 1.067781214143657
$ws.Range("L15").Value = 1.069316219839472
$ws.Range("M15").Value = 1.077961997472902
$ws.Range("N15").Value = 1.069215673814558
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.06258092116269
$ws.Range("D16").Value = 1.065515730290302
$ws.Range("E16").Value = 1.06701347181489
$ws.Range("F16").Value = 1.075808447349418
$ws.Range("I16").Value = 1.044661290009926
$ws.Range("J16").Value = 1.06872009625216
$ws.Range("K16").Value = 1.068851253074073
$ws.Range("L16").Value = 1.070343984481345
$ws.Range("M16").Value = 1.079109848519382
$ws.Range("N16").Value = 1.070237801372747
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.063355197954983
$ws.Range("D17").Value = 1.066257053124135
$ws.Range("E17").Value = 1.06772801241681
$ws.Range("F17").Value = 1.076597931178384
$ws.Range("I17").Value = 1.044867054609844
$ws.Range("J17").Value = 1.069358570760339
$ws.Range("K17").Value = 1.069520730744913
$ws.Range("L17").Value = 1.070986885199202
$ws.Range("M17").Value = 1.079828135370851
$ws.Range("N17").Value = 1.070877182587961
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.063806206299381
$ws.Range("D18").Value = 1.066688884888265
$ws.Range("E18").Value = 1.068144203100292
$ws.Range("F18").Value = 1.077057852401069
$ws.Range("I18").Value = 1.044986627079902
$ws.Range("J18").Value = 1.069730352216909
$ws.Range("K18").Value = 1.069910610630637
$ws.Range("L18").Value = 1.071361241729893
$ws.Range("M18").Value = 1.080246485324614
$ws.Range("N18").Value = 1.071249492016844
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.06395988518857
$ws.Range("D19").Value = 1.066836032610869
$ws.Range("E19").Value = 1.068286014310071
$ws.Range("F19").Value = 1.077214577618607
$ws.Range("I19").Value = 1.04502732273951
$ws.Range("J19").Value = 1.069857013981184
$ws.Range("K19").Value = 1.070043446026596
$ws.Range("L19").Value = 1.071488780416285
$ws.Range("M19").Value = 1.080389028419684
$ws.Range("N19").Value = 1.071376333655352
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.063272189124061
$ws.Range("D20").Value = 1.066177575280642
$ws.Range("E20").Value = 1.067651410063992
$ws.Range("F20").Value = 1.076513286340909
$ws.Range("I20").Value = 1.044845024270044
$ws.Range("J20").Value = 1.069290133821435
$ws.Range("K20").Value = 1.069448965933872
$ws.Range("L20").Value = 1.070917974062935
$ws.Range("M20").Value = 1.07975113375871
$ws.Range("N20").Value = 1.070808648460754
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.061029137862804
$ws.Range("D21").Value = 1.064030114643737
$ws.Range("E21").Value = 1.065581270034274
$ws.Range("F21").Value = 1.074226539944294
$ws.Range("I21").Value = 1.044247085149086
$ws.Range("J21").Value = 1.067439692431381
$ws.Range("K21").Value = 1.067508966807705
$ws.Range("L21").Value = 1.06905468864124
$ws.Range("M21").Value = 1.077669991269411
$ws.Range("N21").Value = 1.068955579231677
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.059613106931606
$ws.Range("D22").Value = 1.062674604854337
$ws.Range("E22").Value = 1.064274201391676
$ws.Range("F22").Value = 1.072783423434226
$ws.Range("I22").Value = 1.043867050168437
$ws.Range("J22").Value = 1.066270401132953
$ws.Range("K22").Value = 1.066283495088562
$ws.Range("L22").Value = 1.06787725773232
$ws.Range("M22").Value = 1.076355768536133
$ws.Range("N22").Value = 1.067784627405487
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.060364325274808
$ws.Range("D23").Value = 1.063393699297693
$ws.Range("E23").Value = 1.06496763260715
$ws.Range("F23").Value = 1.073548963591813
$ws.Range("I23").Value = 1.044068904430359
$ws.Range("J23").Value = 1.066890826980203
$ws.Range("K23").Value = 1.066933691359601
$ws.Range("L23").Value = 1.068502004810394
$ws.Range("M23").Value = 1.077053014850228
$ws.Range("N23").Value = 1.068405934328604
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.063309699138633
$ws.Range("D24").Value = 1.066213489655332
$ws.Range("E24").Value = 1.067686025188331
$ws.Range("F24").Value = 1.076551535463972
$ws.Range("I24").Value = 1.044854980211833
$ws.Range("J24").Value = 1.069321059473958
$ws.Range("K24").Value = 1.06948139526279
$ws.Range("L24").Value = 1.070949114005374
$ws.Range("M24").Value = 1.079785929371013
$ws.Range("N24").Value = 1.070839618031251
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.066703647593408
$ws.Range("D25").Value = 1.069463459578525
$ws.Range("E25").Value = 1.070817564461837
$ws.Range("F25").Value = 1.080013540298184
$ws.Range("I25").Value = 1.045749692131765
$ws.Range("J25").Value = 1.072116582837158
$ws.Range("K25").Value = 1.072413821233331
$ws.Range("L25").Value = 1.073763962715049
$ws.Range("M25").Value = 1.082933308077375
$ws.Range("N25").Value = 1.073639111358278

Write-Output "Applied 380 kV case updates"